$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct BOM line numbering (rows 4-7 renumbered down by 1)
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# Update the saved selection/active cell on the sheet
$ws.Activate()
$ws.Range("B5").Select()
